$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Font.Color = 6710886  # 0x656565 in BGR order used by wdColor? test
$find.Text = ""
$find.Replacement.Font.Color = 255  # red FF0000
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, "", 2, $true, $false, $false, $false)
